$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:E13 (replacing previously internally-rounded values
# with their un-rounded counterparts, per the "Remove internal round and
# float casts" change).
$values = @(
    @(3.657830235493572, 2.002151662949331, 1.570172298926858, 1.378958854892205),
    @(3.555220259392242, 1.971045113597998, 1.558931202156377, 1.372112685789044),
    @(3.607583539008169, 1.993043150492126, 1.57236504845805,  1.383240013895821),
    @(3.534977467695958, 2.020387323550141, 1.57880164550358,  1.399850874782828),
    @(2.968124494335742, 1.875463510908506, 1.548287313870727, 1.384591965316204),
    @(3.230341905891171, 1.961212458634178, 1.587202116579776, 1.390972419996553),
    @(2.048578856136538, 1.497617726625957, 1.347264377506427, 1.307684694103238),
    @(2.83681288685493,  1.812437534362181, 1.524252566768621, 1.359680543457173),
    @(3.499139607109537, 2.033644326230511, 1.602662455506086, 1.403891052508206),
    @(3.48083342614943,  2.025370691254477, 1.599602480699534, 1.400737064552017),
    @(3.018377456737963, 1.892787997103477, 1.56334331424269,  1.373175483619289),
    @(3.337211194759394, 1.989697056952977, 1.587080947628144, 1.392783088423188)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
